$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3229
$ws1.Range("F4").Value = 122
$ws1.Range("F5").Value = 6886
$ws1.Range("F6").Value = 2091
$ws1.Range("F7").Value = 28
$ws1.Range("F8").Value = 76
$ws1.Range("F10").Value = 36
$ws1.Range("F11").Value = 72
$ws1.Range("F14").Value = 189

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3229
$ws4.Range("F3").Value = 19
$ws4.Range("F5").Value = 122
$ws4.Range("F6").Value = 6886
$ws4.Range("F7").Value = 2091
$ws4.Range("F8").Value = 28
$ws4.Range("F9").Value = 76
$ws4.Range("F11").Value = 36
$ws4.Range("F12").Value = 72
$ws4.Range("F15").Value = 189
